$d = $word.ActiveDocument
$d.Content.Find.Execute("30th April 2024", $true, $false, $false, $false, $false, $true, 1, $false, "31st May 2024", 2)
